$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B28 text: "Work on the winning variables ()" -> "Work on the winning/game Over variables ()"
$ws.Range("B28").Value = "Work on the winning/game Over variables ()"

# Remove C28 content: "still using the one from Level One" (clears the now-unused cell entirely)
$ws.Range("C28").Clear()

# Update the active selection to match the new state (B28)
$ws.Range("B28").Select()
